$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / link / non-numeric-looking value updates
$ws.Range("D2").Value = "42.834.13"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.308.71"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("E11").Value = "  +5.46%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").Value = "2.669.84"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").Value = "2.310.50"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "42.762.25"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  -5.08%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("E23").Value = "  +7.78%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("E28").Value = "  +12.68%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("E35").Value = "  -6.49%  "
$ws.Range("E36").Value = "  +1.77%  "
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  +12.43%  "
$ws.Range("D43").Value = "1.975.27"
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("E44").Value = "  +4.85%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").Value = "2.534.35"
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E51").Value = "  +0.67%  "

# Numeric-looking text values in column D must be forced to stay text
# (otherwise Excel auto-converts them to numbers), without leaving a
# lingering custom style on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.786"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0999"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.109"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.49"
$ws.Range("D51").Style = "Normal"
